# Update "想去人数" (number of people interested) counts that changed between
# the previous site generation and the new one (incremented by 1 each).
#
# Sheet "展览" (索引1): F3 1516->1517, F6 13333->13334, F7 13197->13198, F20 257->258
# Sheet "演出" (索引2): F9 30->31
# Sheet "全部类型" (索引4): F4 1516->1517, F8 13333->13334, F9 13197->13198, F27 257->258, F34 30->31

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1517
$wsExhibit.Range("F6").Value = 13334
$wsExhibit.Range("F7").Value = 13198
$wsExhibit.Range("F20").Value = 258

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F9").Value = 31

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1517
$wsAll.Range("F8").Value = 13334
$wsAll.Range("F9").Value = 13198
$wsAll.Range("F27").Value = 258
$wsAll.Range("F34").Value = 31
